$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 969.89655
$ws.Range("J17").Value = 1096.826
$ws.Range("L17").Value = 3290.478
$ws.Range("N17").Value = -3626.478

# Sheet ALC, Row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3072.6128
$ws.Range("I64").Value = 2920.6
$ws.Range("J64").Value = 3145
$ws.Range("K64").Value = 2920.6
$ws.Range("L64").Value = 3145
$ws.Range("M64").Value = -2672.6
$ws.Range("N64").Value = -3641

# Sheet ALC, Row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3072.6128
$ws.Range("I67").Value = 2920.6
$ws.Range("J67").Value = 3145
$ws.Range("K67").Value = 2920.6
$ws.Range("L67").Value = 3145
$ws.Range("M67").Value = -2062.6
$ws.Range("N67").Value = -4861

# Sheet ALC, Row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3416.4285
$ws.Range("I74").Value = 3053.6365
$ws.Range("K74").Value = 3053.6365
$ws.Range("M74").Value = -2117.6365

# Sheet ALC, Row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3413.2856
$ws.Range("I76").Value = 2899.125
$ws.Range("J76").Value = 4098.8335
$ws.Range("K76").Value = 2899.125
$ws.Range("L76").Value = 4098.8335
$ws.Range("M76").Value = -2584.125
$ws.Range("N76").Value = -4728.8335

# Sheet ALC, Row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3416.4285
$ws.Range("I77").Value = 3053.6365
$ws.Range("K77").Value = 15268.1825
$ws.Range("M77").Value = -10588.1825

# Sheet ALC, Row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3413.2856
$ws.Range("I79").Value = 2899.125
$ws.Range("J79").Value = 4098.8335
$ws.Range("K79").Value = 2899.125
$ws.Range("L79").Value = 4098.8335
$ws.Range("M79").Value = -1807.125
$ws.Range("N79").Value = -6282.8335

# Sheet ALC, Row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1932.7273
$ws.Range("I129").Value = 657.3333
$ws.Range("J129").Value = 2815.6924
$ws.Range("K129").Value = 1971.9999
$ws.Range("L129").Value = 8447.0772
$ws.Range("M129").Value = 3028.0001
$ws.Range("N129").Value = -18447.0772

# Sheet ALC, Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2384.395
$ws.Range("J138").Value = 2928.6086
$ws.Range("L138").Value = 8785.825800000001
$ws.Range("N138").Value = -19065.8258

# Sheet ALC, Row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7733
$ws.Range("I141").Value = 1859.6875
$ws.Range("J141").Value = 31226.25
$ws.Range("K141").Value = 5579.0625
$ws.Range("L141").Value = 93678.75
$ws.Range("M141").Value = -399.0625
$ws.Range("N141").Value = -104038.75

# Sheet ARM, Row 30
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 950
$ws.Range("I30").Value = 900
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 900
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = -750
$ws.Range("N30").Value = -1300

# Sheet ARM, Row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1020.4
$ws.Range("I88").Value = 1025.5
$ws.Range("J88").Value = 1000
$ws.Range("K88").Value = 1025.5
$ws.Range("L88").Value = 1000
$ws.Range("M88").Value = -619.5
$ws.Range("N88").Value = -1812

# Sheet ARM, Row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1020.4
$ws.Range("I91").Value = 1025.5
$ws.Range("J91").Value = 1000
$ws.Range("K91").Value = 1025.5
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = 378.5
$ws.Range("N91").Value = -3808

# Sheet CRP, Row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2000301.4
$ws.Range("I6").Value = 2500176.8
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 2500176.8
$ws.Range("L6").Value = 800
$ws.Range("M6").Value = -2500063.8
$ws.Range("N6").Value = -1026

# Sheet CRP, Row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 465.0909
$ws.Range("I19").Value = 471.6
$ws.Range("K19").Value = 471.6
$ws.Range("M19").Value = -301.6

# Sheet CRP, Row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 465.0909
$ws.Range("I24").Value = 471.6
$ws.Range("K24").Value = 471.6
$ws.Range("M24").Value = -301.6

# Sheet CRP, Row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 15213.75
$ws.Range("I62").Value = 2916.6667
$ws.Range("J62").Value = 22592
$ws.Range("K62").Value = 2916.6667
$ws.Range("L62").Value = 22592
$ws.Range("M62").Value = -2292.6667
$ws.Range("N62").Value = -23840

# Sheet CRP, Row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 15213.75
$ws.Range("I65").Value = 2916.6667
$ws.Range("J65").Value = 22592
$ws.Range("K65").Value = 14583.3335
$ws.Range("L65").Value = 112960
$ws.Range("M65").Value = -11463.3335
$ws.Range("N65").Value = -119200

# Sheet CRP, Row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 32261756
$ws.Range("I86").Value = 41669804
$ws.Range("J86").Value = 5598.2856
$ws.Range("K86").Value = 41669804
$ws.Range("L86").Value = 5598.2856
$ws.Range("M86").Value = -41668681
$ws.Range("N86").Value = -7844.2856

# Sheet CRP, Row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 32261756
$ws.Range("I89").Value = 41669804
$ws.Range("J89").Value = 5598.2856
$ws.Range("K89").Value = 208349020
$ws.Range("L89").Value = 27991.428
$ws.Range("M89").Value = -208343404
$ws.Range("N89").Value = -39223.428

# Sheet CRP, Row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 589868.9
$ws.Range("I105").Value = 668277.3
$ws.Range("J105").Value = 1805.5
$ws.Range("K105").Value = 668277.3
$ws.Range("L105").Value = 1805.5
$ws.Range("M105").Value = -666530.3
$ws.Range("N105").Value = -5299.5

# Sheet CUL, Row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2862.8
$ws.Range("I22").Value = 2540.4
$ws.Range("J22").Value = 3024
$ws.Range("K22").Value = 7621.200000000001
$ws.Range("L22").Value = 9072
$ws.Range("M22").Value = -7452.200000000001
$ws.Range("N22").Value = -9410

# Sheet CUL, Row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 2862.8
$ws.Range("I27").Value = 2540.4
$ws.Range("J27").Value = 3024
$ws.Range("K27").Value = 7621.200000000001
$ws.Range("L27").Value = 9072
$ws.Range("M27").Value = -7519.200000000001
$ws.Range("N27").Value = -9276

# Sheet CUL, Row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 19566.363
$ws.Range("J34").Value = 30707.143
$ws.Range("L34").Value = 92121.429
$ws.Range("N34").Value = -92289.429

# Sheet CUL, Row 102
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 5010.857
$ws.Range("I102").Value = 1999.5
$ws.Range("J102").Value = 6215.4
$ws.Range("K102").Value = 5998.5
$ws.Range("L102").Value = 18646.2
$ws.Range("M102").Value = -3564.5
$ws.Range("N102").Value = -23514.2

# Sheet GSM, Row 10
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 300
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# Sheet GSM, Row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7430.857
$ws.Range("I70").Value = 8028
$ws.Range("K70").Value = 8028
$ws.Range("M70").Value = -7758

# Sheet GSM, Row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7430.857
$ws.Range("I73").Value = 8028
$ws.Range("K73").Value = 8028
$ws.Range("M73").Value = -7092

# Sheet GSM, Row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2450.5454
$ws.Range("I80").Value = 2250
$ws.Range("K80").Value = 2250
$ws.Range("M80").Value = -1252

# Sheet GSM, Row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2450.5454
$ws.Range("I83").Value = 2250
$ws.Range("K83").Value = 11250
$ws.Range("M83").Value = -6258

# Sheet GSM, Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2880.7104
$ws.Range("I132").Value = 1681.4736
$ws.Range("J132").Value = 4079.9473
$ws.Range("K132").Value = 5044.4208
$ws.Range("L132").Value = 12239.8419
$ws.Range("M132").Value = -2514.4208
$ws.Range("N132").Value = -17299.8419

# Sheet LTW, Row 120
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H120").Value = 30739.4
$ws.Range("J120").Value = 28249.75
$ws.Range("L120").Value = 28249.75
$ws.Range("N120").Value = -37925.75

# Sheet WVR, Row 9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 10000000
$ws.Range("I9").Value = 10000000
$ws.Range("K9").Value = 10000000
$ws.Range("M9").Value = -9999860
